# New Varying Gaussian data for Spheres, Cones and Propellor
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the static correction-factor values in C7, C15 and C23 with the
# Varying Gaussian formula that derives them from the refractive index (B column).
$ws.Range("C7").Formula  = "=((((1.4435*(B7*B7))/((B7*B7)-0.020216))+1)^0.5)/1.328"
$ws.Range("C15").Formula = "=((((1.4435*(B15*B15))/((B15*B15)-0.020216))+1)^0.5)/1.328"
$ws.Range("C23").Formula = "=((((1.4435*(B23*B23))/((B23*B23)-0.020216))+1)^0.5)/1.328"

# Update the saved cursor / selection position on the sheet.
$ws.Range("E24").Select()
